$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("purchase_page")

# --- Add a new "Test Case Summary" block (rows 22-28), mirroring the
#     existing blocks at rows 2-8 and 11-17, for the new CrabVPN test run. ---

# Copy the cell formatting (fonts, fills, borders, number formats) from the
# most recent existing block (rows 11-17) down onto the new block's rows.
# (Done as two pieces - D11:E15 and D17:E17 - so the blank spacer row
# between them, row 16, doesn't create an empty row 27 in the new block.)
$fmtSrc1 = $ws.Range("D11:E15")
$fmtSrc1.Copy()
$ws.Range("D22:E26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$fmtSrc2 = $ws.Range("D17:E17")
$fmtSrc2.Copy()
$ws.Range("D28:E28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row for the new block.
$ws.Range("D22").Value2 = " Test Case Summary (06-04-24)"
$ws.Range("F22").Value2 = "(CrabVPN)"

# Summary counts for the new run.
$ws.Range("D23").Value2 = "Execute"
$ws.Range("E23").Value2 = 33

$ws.Range("D24").Value2 = "Pass"
$ws.Range("E24").Value2 = 30

$ws.Range("D25").Value2 = "Fail"
$ws.Range("E25").Value2 = 0

$ws.Range("D26").Value2 = "Not Tested"
$ws.Range("E26").Value2 = 3

# Row heights for the new header/summary rows match the existing blocks.
$ws.Rows.Item(22).RowHeight = 18
$ws.Rows.Item(23).RowHeight = 18
$ws.Rows.Item(24).RowHeight = 18
$ws.Rows.Item(25).RowHeight = 18
$ws.Rows.Item(26).RowHeight = 18

# Merge the header cells for the new block, like the other header rows.
$ws.Range("D22:E22").Merge()

# Keep focus/selection on the newly added block.
$ws.Range("R7").Select()
